$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.322.45"
$ws.Range("E2").Value = "  -1.08%  "

$ws.Range("D3").Value = "3.899.84"
$ws.Range("E3").Value = "  -0.74%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  -0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "485.45"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.56%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "146.08"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.39%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +0.18%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -0.01%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.741"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +2.64%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.177"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.19%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0000355"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -0.49%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.87"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.70%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.49"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.09%  "

$ws.Range("D14").Value = "4.512.35"
$ws.Range("E14").Value = "  -1.06%  "

$ws.Range("D15").Value = "3.922.15"
$ws.Range("E15").Value = "  -0.92%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "14.28"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -1.88%  "

$ws.Range("E17").Value = "  -0.61%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "19.98"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.14"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +0.71%  "

$ws.Range("D20").Value = "68.319.46"
$ws.Range("E20").Value = "  -1.14%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "430.64"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.85%  "

$ws.Range("E22").Value = "  +7.38%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "14.76"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +1.06%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "12.45"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +19.37%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "88.99"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +1.12%  "

$ws.Range("E26").Value = "  +4.34%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.02"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -5.32%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "37.23"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.36%  "

$ws.Range("E29").Value = "  -3.34%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "720.43"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.97%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "13.48"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +2.00%  "

$ws.Range("E32").Value = "  +0.44%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.92"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.69%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "61.78"
$ws.Range("D34").ClearFormats()

$ws.Range("B35").Value = "PEPE"
$ws.Range("C35").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D35").Value = "0.0₃0876"
$ws.Range("E35").Value = "  -6.47%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.04"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +8.81%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "40.92"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -0.48%  "

$ws.Range("E38").Value = "  +18.16%  "

$ws.Range("B39").Value = "Kaspa"
$ws.Range("C39").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.146"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -3.87%  "

$ws.Range("B40").Value = "Dai"
$ws.Range("C40").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.00"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +0.14%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "3.00"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +9.17%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0497"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +5.57%  "

$ws.Range("E43").Value = "  +3.11%  "

$ws.Range("E45").Value = "  +1.32%  "

$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").Value = "0.0₆0366"
$ws.Range("E46").Value = "  +32.47%  "

$ws.Range("B47").Value = "FirstDigitalUSD"
$ws.Range("C47").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.00"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -0.13%  "

$ws.Range("E48").Value = "  +5.73%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.37"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -1.14%  "

$ws.Range("E50").Value = "  -2.22%  "

$ws.Range("E51").Value = "  -2.67%  "
